$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. "Date created" (B3): was a text placeholder ("Unknown"); now the actual
#    year the dataset was created, stored as a number, left/top aligned and
#    wrapped.
# ---------------------------------------------------------------------------
$b3 = $ws.Cells.Item(3, 2)
$b3.Value = 2007
$b3.HorizontalAlignment = -4131   # xlLeft
$b3.VerticalAlignment = -4160     # xlTop
$b3.WrapText = $true

# ---------------------------------------------------------------------------
# 2. "Taxa" attribute description (B10) gets expanded wording.
# ---------------------------------------------------------------------------
$ws.Cells.Item(10, 2).Value = "Target vertebrate taxa for housing. This may later include insects (e.g. Hymenopterans, bees) or other artificial wildlife housing. "

# ---------------------------------------------------------------------------
# 3. Four new attribute rows (RoostSp, Year, Month, Day) are appended below
#    the existing "Taxa" row, matching its formatting.
# ---------------------------------------------------------------------------
$ws.Range("A10:F10").Copy()
$ws.Range("A11:F14").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Rows.Item(11).RowHeight = 24.75
$ws.Rows.Item(12).RowHeight = 72.75
$ws.Rows.Item(13).RowHeight = 72.75
$ws.Rows.Item(14).RowHeight = 84.75

# --- Row 11: RoostSp ---------------------------------------------------
$ws.Cells.Item(11, 1).Value = "RoostSp"
$ws.Cells.Item(11, 2).Value = "Species that the wildlfe box was attached (E.g. birdbox attached to a large Douglas Fir tree). "
$ws.Cells.Item(11, 4).Value = "String"

$e11 = $ws.Cells.Item(11, 5)
$e11rest = " `nE.g. Pseudotsuga menziesii."
$e11.Value = "Values:" + $e11rest
$e11.Font.Name = "Times New Roman"
$e11.Font.Size = 9
$e11.Characters(1, 7).Font.Underline = $true
$e11.Characters(8, $e11rest.Length).Font.Underline = $false

# --- Row 12: Year --------------------------------------------------------
$ws.Cells.Item(12, 1).Value = "Year"
$ws.Cells.Item(12, 2).Value = "Year the data was recorded."
$ws.Cells.Item(12, 3).Value = "Date"
$ws.Cells.Item(12, 4).Value = "String"

$e12 = $ws.Cells.Item(12, 5)
$e12rest = " `nyyyy. E.g. 2008. `nNULL = neither the original meta-data nor accompanying report provided the year of creation. "
$e12.Value = "Values:" + $e12rest
$e12.Font.Name = "Times New Roman"
$e12.Font.Size = 9
$e12.Characters(1, 7).Font.Underline = $true
$e12.Characters(8, $e12rest.Length).Font.Underline = $false

# --- Row 13: Month ---------------------------------------------------------
$ws.Cells.Item(13, 1).Value = "Month"
$ws.Cells.Item(13, 2).Value = "Month the data was recorded."
$ws.Cells.Item(13, 3).Value = "Date"
$ws.Cells.Item(13, 4).Value = "String"

$e13 = $ws.Cells.Item(13, 5)
$e13rest = " `n1-12. E.g. 2=February.`nNULL = neither the original meta-data nor accompanying report provided the month of creation. "
$e13.Value = "Values:" + $e13rest
$e13.Font.Name = "Times New Roman"
$e13.Font.Size = 9
$e13.Characters(1, 7).Font.Underline = $true
$e13.Characters(8, $e13rest.Length).Font.Underline = $false

# --- Row 14: Day -----------------------------------------------------------
$ws.Cells.Item(14, 1).Value = "Day"
$ws.Cells.Item(14, 2).Value = "Day the data was recorded."
$ws.Cells.Item(14, 3).Value = "Date"
$ws.Cells.Item(14, 4).Value = "String"

$e14 = $ws.Cells.Item(14, 5)
$e14rest = " `n1-31. E.g. 15=the 15th day of a month. `nNULL = neither the original meta-data nor accompanying report provided the day of creation."
$e14.Value = "Values:" + $e14rest
$e14.Font.Name = "Times New Roman"
$e14.Font.Size = 9
$e14.Characters(1, 7).Font.Underline = $true
$e14.Characters(8, $e14rest.Length).Font.Underline = $false

# ---------------------------------------------------------------------------
# 4. Update the active selection to the newly added rows, like the author
#    left it selected after editing.
# ---------------------------------------------------------------------------
$ws.Rows.Item("12:14").Select()
